# submissions_daily_matrix.xlsx — add a new daily column "2026-02-13"
# This inserts a new column before the existing "total_files" column (F),
# shifting total_files -> G and unique_days -> H automatically, then
# populates the new date column's header and data (all zero, since no
# submissions have been recorded yet for that day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; existing F (total_files) becomes G,
# existing G (unique_days) becomes H.
$ws.Columns("F:F").Insert()

# --- Header cell F1: new date column, styled like the other date headers (D1/E1) ---
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "2026-02-13"

# Copy the exact formatting (font/alignment/fill/number format) from E1 onto F1
# so it matches the other date-header cells precisely.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data cells F2:F109: new date column has no submissions yet -> all 0 ---
$ws.Range("F2:F109").Value = 0

# --- Column width: new date column matches the other date columns (width 12) ---
# (ColumnWidth uses character units that differ from the stored sheet "width" by
#  a constant offset of ~0.8333 on this font; subtract it so the saved width is 12.)
$ws.Columns("F:F").ColumnWidth = 11.166666666666666

Write-Output "Applied: inserted 2026-02-13 column before total_files/unique_days."
